# PROS-12599 - MARSRU - KPIs PSS 2020
#
# Adds 4 new KPI rows (product codes 5070, 5032, 5033, 5034) to both the
# "KPI with Codes" and "KPI with Names" sheets, preceded by a blank
# separator row (row 152, matching the same blank-row convention already
# used elsewhere in the sheet, e.g. row 144), using the same
# layout/formula conventions as the existing rows.
#
# NOTE: rows are written directly at their target row numbers (152 and
# 158-161) rather than via an Insert/shift, because in the source sheet
# row 152 (and 157) were already "missing" row numbers (a gap between 151
# and 153, and between 156 and 158) - the edit fills that gap in and
# appends after row 156 without moving any existing row.

$wb = $excel.ActiveWorkbook

function Add-BlankSeparatorRow($ws, $rowNum, $templateRow) {
    # Mirror the look of the existing blank separator row (e.g. row 144):
    # same formatting, but completely empty cells.
    $src = $ws.Range("A" + $templateRow + ":K" + $templateRow)
    $dst = $ws.Range("A" + $rowNum + ":K" + $rowNum)
    $src.Copy($dst)
    $dst.ClearContents()
    # Touch a benign formatting property so the (otherwise fully empty)
    # row actually materialises in the sheet.
    $dst.WrapText = $false
}

function Add-CodeRow($ws, $row, $aVal, $aIsText, $bVal, $cVal, $eVal, $styleTemplateRow, $kFormula) {
    # Copy an existing, fully-formatted data row as a template: this
    # brings over correct formatting for every column, AND the F:J helper
    # cells (which hold literal fragments of SQL-ish text, several of
    # which start with a leading "'" that a plain .Value assignment would
    # otherwise strip as a "quote prefix").
    $srcRange = $ws.Range("A" + $styleTemplateRow + ":K" + $styleTemplateRow)
    $dstRange = $ws.Range("A" + $row + ":K" + $row)
    $srcRange.Copy($dstRange)

    $ws.Rows.Item($row).RowHeight = 12.85

    $aCell = $ws.Range("A" + $row)
    if ($aIsText) {
        # Forces a numeric-looking literal ("5032" etc.) to be stored as
        # text instead of being auto-coerced to a number.
        $aCell.NumberFormat = "@"
        $aCell.Value = $aVal
        $aCell.NumberFormat = "General"
    } else {
        $aCell.Value = $aVal
    }

    $ws.Range("B" + $row).Value = $bVal
    $ws.Range("C" + $row).Value = $cVal

    $dFormula = '=IF(EXACT(C' + $row + ',"Boolean"),"''10''","NULL")'
    $ws.Range("D" + $row).Formula = $dFormula

    $ws.Range("E" + $row).Value = $eVal

    # F:J are left as copied from the template row (identical literal
    # text fragments in every data row).

    $ws.Range("K" + $row).Formula = $kFormula
}

$wsCodes = $wb.Worksheets.Item("KPI with Codes")
$wsNames = $wb.Worksheets.Item("KPI with Names")

$kFormulaCodesTpl = '=CONCATENATE(F{0},E{0},G{0},A{0},H{0},A{0},I{0},D{0},J{0})'
$kFormulaNamesTpl = '=CONCATENATE(F{0},E{0},G{0},CONCATENATE(A{0},"-RUS"),H{0},CONCATENATE(A{0}," - ",B{0}),I{0},D{0},J{0})'

$bVal158 = "Секц лаком д/кош прим к корм д/кош и лаком д/соб к корм д/соб. Допуск разм между влаж и сух"
$bVal159 = "PERFECT FIT вз.кош чувст инд 24*85г [Кол-во горизонтальных фэйсов на основной полке]"
$bVal160 = "Sheba Плежер куриц и кролик 24*85г [Кол-во горизонтальных фэйсов на основной полке]"
$bVal161 = "Cesar гов. с овощами 100г [Кол-во горизонтальных фэйсов на основной полке]"

foreach ($ws in @($wsCodes, $wsNames)) {
    Add-BlankSeparatorRow $ws 152 144

    if ($ws.Name -eq "KPI with Codes") {
        $eBase = 1441
        $kTpl = $kFormulaCodesTpl
    } else {
        $eBase = 2739
        $kTpl = $kFormulaNamesTpl
    }

    $e158 = $eBase
    $e159 = $eBase + 1
    $e160 = $eBase + 2
    $e161 = $eBase + 3

    $k158 = $kTpl -f 158
    $k159 = $kTpl -f 159
    $k160 = $kTpl -f 160
    $k161 = $kTpl -f 161

    Add-CodeRow $ws 158 5070 $false $bVal158 "Boolean" $e158 150 $k158
    Add-CodeRow $ws 159 "5032" $true $bVal159 "Int" $e159 150 $k159
    Add-CodeRow $ws 160 "5033" $true $bVal160 "Int" $e160 150 $k160
    Add-CodeRow $ws 161 "5034" $true $bVal161 "Int" $e161 150 $k161
}

$wsNames.Range("K158:K161").Select()
$wsCodes.Activate()
$wsCodes.Range("K158:K161").Select()
